# Apply updated time-to-discovery simulation values to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C2" = 294;    "D2" = 293.5;
    "C3" = 61;     "D3" = 52;
    "B4" = 1673;   "C4" = 1693;   "D4" = 1683;
    "C5" = 34;     "D5" = 32.5;
    "C6" = 262;    "D6" = 243;
    "C7" = 27;     "D7" = 57;
    "C8" = 108;    "D8" = 60.5;
    "C9" = 62;     "D9" = 53;
    "C10" = 270;   "D10" = 245;
    "C11" = 109;   "D11" = 61.5;
    "C12" = 60;    "D12" = 47;
    "C13" = 185;   "D13" = 159;
    "C14" = 126;   "D14" = 113.5;
    "C15" = 35;    "D15" = 33.5;
    "C16" = 103;   "D16" = 94.5;
    "C17" = 37;    "D17" = 37;
    "C18" = 135;   "D18" = 113.5;
    "C19" = 131;   "D19" = 123;
    "C20" = 174;   "D20" = 150;
    "C21" = 53;    "D21" = 30.5;
    "C23" = 57;    "D23" = 41;
    "C24" = 241;   "D24" = 243;
    "C25" = 70;    "D25" = 58;
    "C26" = 41;    "D26" = 53.5;
    "C27" = 541;   "D27" = 509.5;
    "B28" = 169.56; "C28" = 190.5769230769231;
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
